# Swap the presentation's theme colour scheme between the custom "Integral"
# palette and the default "Office" palette (Design > Variants > Colors).
#
# The deck currently uses the "Integral" theme colours for its slide master
# (ppt/theme/theme1.xml) while its notes master (ppt/theme/theme2.xml) still
# carries the original default "Office" colours. The edit applies the
# default "Office" colour scheme to the deck's theme, matching the colour
# values that were previously only used by the notes master.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = 0x00000000   # Dark 1    -> 000000
$colors.Item(2).RGB  = 0x00FFFFFF   # Light 1   -> FFFFFF
$colors.Item(3).RGB  = 0x006A5444   # Dark 2    -> 44546A
$colors.Item(4).RGB  = 0x00E6E6E7   # Light 2   -> E7E6E6
$colors.Item(5).RGB  = 0x00D59B5B   # Accent 1  -> 5B9BD5
$colors.Item(6).RGB  = 0x00317DED   # Accent 2  -> ED7D31
$colors.Item(7).RGB  = 0x00A5A5A5   # Accent 3  -> A5A5A5
$colors.Item(8).RGB  = 0x0000C0FF   # Accent 4  -> FFC000
$colors.Item(9).RGB  = 0x00C47244   # Accent 5  -> 4472C4
$colors.Item(10).RGB = 0x0047AD70   # Accent 6  -> 70AD47
$colors.Item(11).RGB = 0x00C16305   # Hyperlink -> 0563C1
$colors.Item(12).RGB = 0x00724F95   # Followed Hyperlink -> 954F72
